$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting the old Weight/Deflection/Stress
# columns (C:E) one column to the right (D:F).
$ws.Columns("C:C").Insert()

# --- Header row ---------------------------------------------------
$ws.Range("A1").Value = "Half Span"
$ws.Range("B1").Value = "Root Chord"
$ws.Range("C1").Value = "Tip Chord"
# D1/E1/F1 (Weight/Deflection/Stress) keep their original text - no change needed.

# --- Data rows ------------------------------------------------------
# Column A used to hold "Span" sample values (3,4,3); it now holds what
# used to be the "Chord Length" values (30,35,40).
$ws.Range("A2").Value = 30
$ws.Range("A3").Value = 35
$ws.Range("A4").Value = 40

# Column B now holds new "Root Chord" values.
$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 9

# Column C is the new "Tip Chord" formula column (= Root Chord * 0.45).
$ws.Range("C2").Formula = "=B2*0.45"
$ws.Range("C3:C4").Formula = "=B3*0.45"

# --- Sheet view tweaks ----------------------------------------------
[void]$ws.Range("A2").Select()
